$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.022.47"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "1.576.97"
$ws.Range("E3").Value = "  -2.01%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").Value = "299.13"
$ws.Range("E6").Value = "  -1.29%  "

$ws.Range("D7").Value = "0.3747"
$ws.Range("E7").Value = "  -0.86%  "

$ws.Range("D8").Value = "0.3569"
$ws.Range("E8").Value = "  -2.85%  "

$ws.Range("D9").Value = "50.06"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").Value = "1.002"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").Value = "1.212"
$ws.Range("E11").Value = "  -5.04%  "

$ws.Range("D12").Value = "0.07954"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").Value = "21.81"
$ws.Range("E13").Value = "  -5.84%  "

$ws.Range("D14").Value = "6.439"
$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("D15").Value = "7.274"
$ws.Range("E15").Value = "  -4.45%  "

$ws.Range("D16").Value = "0.00001217"
$ws.Range("E16").Value = "  -4.29%  "

$ws.Range("D17").Value = "1.581.21"
$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").Value = "91.68"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").Value = "0.06749"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  -3.76%  "

$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").Value = "6.358"
$ws.Range("E22").Value = "  -3.66%  "

$ws.Range("D23").Value = "23.005.96"
$ws.Range("E23").Value = "  -0.88%  "

$ws.Range("D24").Value = "12.67"
$ws.Range("E24").Value = "  -3.31%  "

$ws.Range("D25").Value = "2.364"
$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("D26").Value = "2.794"
$ws.Range("E26").Value = "  -4.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.60"
$ws.Range("E27").Value = "  -2.47%  "

$ws.Range("D28").Value = "147.05"
$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("D29").Value = "5.174"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").Value = "131.52"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").Value = "2.321"
$ws.Range("E31").Value = "  -3.95%  "

$ws.Range("D32").Value = "6.461"
$ws.Range("E32").Value = "  -7.78%  "

$ws.Range("D33").Value = "1.758.26"
$ws.Range("E33").Value = "  -1.56%  "

$ws.Range("D34").Value = "0.9263"
$ws.Range("E34").Value = "  -5.71%  "

$ws.Range("D35").Value = "0.07303"
$ws.Range("E35").Value = "  -5.88%  "

$ws.Range("D36").Value = "0.02655"
$ws.Range("E36").Value = "  -4.99%  "

$ws.Range("D37").Value = "0.2482"
$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("D38").Value = "0.08726"
$ws.Range("E38").Value = "  -1.67%  "

$ws.Range("D39").Value = "9.828"
$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("D40").Value = "5.944"
$ws.Range("E40").Value = "  -5.84%  "

$ws.Range("D41").Value = "1.337"
$ws.Range("E41").Value = "  -4.38%  "

$ws.Range("D42").Value = "0.6826"
$ws.Range("E42").Value = "  -4.95%  "

$ws.Range("D43").Value = "11.75"
$ws.Range("E43").Value = "  -8.25%  "

$ws.Range("D44").Value = "14.67"
$ws.Range("E44").Value = "  -8.95%  "

$ws.Range("D45").Value = "0.6325"
$ws.Range("E45").Value = "  -4.72%  "

$ws.Range("D46").Value = "3.963"
$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").Value = "2.239"
$ws.Range("E47").Value = "  -2.93%  "

$ws.Range("D48").Value = "130.85"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("E49").Value = "  -2.33%  "

$ws.Range("D50").Value = "1.178"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").Value = "1.165"
$ws.Range("E51").Value = "  -1.90%  "
